# Staging.Gender template regenerated after Meerkat DB changes:
# a new "BusinessKey" column is inserted as the first data column on the
# header row, pushing the existing Code / Gender_ID / Name headers one
# column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the current header text before it gets shifted.
$oldCode     = $ws.Range("A2").Value2
$oldGenderId = $ws.Range("B2").Value2
$oldName     = $ws.Range("C2").Value2

# Shift the existing headers one column to the right.
$ws.Range("D2").Value = $oldName
$ws.Range("C2").Value = $oldGenderId
$ws.Range("B2").Value = $oldCode

# New first header column.
$ws.Range("A2").Value = "BusinessKey"

# Give the new header cell the same bold + underline look as the rest of
# the header row (re-uses the existing header style rather than creating
# a new one).
$ws.Range("D2").Font.Bold = $true
$ws.Range("D2").Font.Underline = $true
